# Automatic update of files.
# Rewrites rows 53-62 on the "Artfynd" sheet: species records were
# reshuffled among the rows (site coordinates Q/R stay tied to a given
# row while the A/B/D/E/F/G/H species-identifying fields move between
# rows), plus row 56 gains a "larv/nymf" age/stage annotation (columns
# J/K/L/M/N/AF) that row 59 loses.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 53 ---------------------------------------------------------
$ws.Cells.Item(53, 1).Value = 112230614   # A53
$ws.Cells.Item(53, 2).Value = 78647       # B53
$ws.Cells.Item(53, 17).Value = 571792     # Q53
$ws.Cells.Item(53, 18).Value = 6697651    # R53

# --- Row 54 ---------------------------------------------------------
$ws.Cells.Item(54, 1).Value = 112230603   # A54
$ws.Cells.Item(54, 2).Value = 78647       # B54
$ws.Cells.Item(54, 17).Value = 572018     # Q54
$ws.Cells.Item(54, 18).Value = 6697738    # R54

# --- Row 55 ---------------------------------------------------------
$ws.Cells.Item(55, 1).Value = 112230604              # A55
$ws.Cells.Item(55, 2).Value = 102192                 # B55
$ws.Cells.Item(55, 5).Value = 222412                 # E55
$ws.Cells.Item(55, 6).Value = "Tibast"                # F55
$ws.Cells.Item(55, 7).Value = "Daphne mezereum"       # G55
$ws.Cells.Item(55, 8).Value = "L."                    # H55
$ws.Cells.Item(55, 17).Value = 571996                # Q55
$ws.Cells.Item(55, 18).Value = 6697876               # R55

# --- Row 56 (gains J/K/L/M/N/AF; K56 = "larv/nymf") -----------------
$ws.Cells.Item(56, 1).Value = 112230612              # A56
$ws.Cells.Item(56, 2).Value = 12274                  # B56
$ws.Cells.Item(56, 5).Value = 102016                 # E56
$ws.Cells.Item(56, 6).Value = "Gropig brunbagge"      # F56
$ws.Cells.Item(56, 7).Value = "Zilora ferruginea"     # G56
$ws.Cells.Item(56, 8).Value = "(Paykull, 1798)"       # H56
$ws.Cells.Item(56, 10).Value = ""                    # J56 (new, blank)
$ws.Cells.Item(56, 10).Style = "Normal"
$ws.Cells.Item(56, 11).Value = "larv/nymf"            # K56 (new)
$ws.Cells.Item(56, 12).Value = ""                    # L56 (new, blank)
$ws.Cells.Item(56, 12).Style = "Normal"
$ws.Cells.Item(56, 13).Value = ""                    # M56 (new, blank)
$ws.Cells.Item(56, 13).Style = "Normal"
$ws.Cells.Item(56, 14).Value = ""                    # N56 (new, blank)
$ws.Cells.Item(56, 14).Style = "Normal"
$ws.Cells.Item(56, 17).Value = 571800                # Q56
$ws.Cells.Item(56, 18).Value = 6697623               # R56
$ws.Cells.Item(56, 32).Value = ""                    # AF56 (new, blank)
$ws.Cells.Item(56, 32).Style = "Normal"

# --- Row 57 ---------------------------------------------------------
$ws.Cells.Item(57, 1).Value = 112230605              # A57
$ws.Cells.Item(57, 2).Value = 99874                  # B57
$ws.Cells.Item(57, 5).Value = 221235                 # E57
$ws.Cells.Item(57, 6).Value = "Vårärt"                # F57
$ws.Cells.Item(57, 7).Value = "Lathyrus vernus"       # G57
$ws.Cells.Item(57, 8).Value = "(L.) Bernh."           # H57
$ws.Cells.Item(57, 17).Value = 571995                # Q57

# --- Row 58 ---------------------------------------------------------
$ws.Cells.Item(58, 1).Value = 112230608              # A58
$ws.Cells.Item(58, 2).Value = 99874                  # B58
$ws.Cells.Item(58, 4).Value = "LC"                    # D58
$ws.Cells.Item(58, 5).Value = 221235                 # E58
$ws.Cells.Item(58, 6).Value = "Vårärt"                # F58
$ws.Cells.Item(58, 7).Value = "Lathyrus vernus"       # G58
$ws.Cells.Item(58, 8).Value = "(L.) Bernh."           # H58
$ws.Cells.Item(58, 17).Value = 571931                # Q58
$ws.Cells.Item(58, 18).Value = 6697694               # R58

# --- Row 59 (loses J/K/L/M/N/AF) -------------------------------------
$ws.Cells.Item(59, 1).Value = 112230610              # A59
$ws.Cells.Item(59, 2).Value = 90480                  # B59
$ws.Cells.Item(59, 4).Value = "LC"                    # D59
$ws.Cells.Item(59, 5).Value = 4769                   # E59
$ws.Cells.Item(59, 6).Value = "Svavelriska"           # F59
$ws.Cells.Item(59, 7).Value = "Lactarius scrobiculatus" # G59
$ws.Cells.Item(59, 8).Value = "(Scop.:Fr.) Fr."       # H59
$ws.Cells.Item(59, 10).ClearContents()               # J59 removed
$ws.Cells.Item(59, 11).ClearContents()               # K59 removed ("larv/nymf")
$ws.Cells.Item(59, 12).ClearContents()               # L59 removed
$ws.Cells.Item(59, 13).ClearContents()               # M59 removed
$ws.Cells.Item(59, 14).ClearContents()               # N59 removed
$ws.Cells.Item(59, 17).Value = 571853                # Q59
$ws.Cells.Item(59, 18).Value = 6697760               # R59
$ws.Cells.Item(59, 32).ClearContents()               # AF59 removed

# --- Row 60 ---------------------------------------------------------
$ws.Cells.Item(60, 1).Value = 112230606              # A60
$ws.Cells.Item(60, 2).Value = 56575                  # B60
$ws.Cells.Item(60, 4).Value = "NT"                    # D60
$ws.Cells.Item(60, 5).Value = 103021                 # E60
$ws.Cells.Item(60, 6).Value = "Talltita"              # F60
$ws.Cells.Item(60, 7).Value = "Poecile montanus"      # G60
$ws.Cells.Item(60, 8).Value = "(Conrad von Baldenstein, 1827)" # H60
$ws.Cells.Item(60, 17).Value = 571961                # Q60
$ws.Cells.Item(60, 18).Value = 6697705               # R60

# --- Row 61 ---------------------------------------------------------
$ws.Cells.Item(61, 1).Value = 112230611              # A61
$ws.Cells.Item(61, 2).Value = 4711                   # B61
$ws.Cells.Item(61, 5).Value = 100299                 # E61
$ws.Cells.Item(61, 6).Value = "Thomsons trägnagare"   # F61
$ws.Cells.Item(61, 7).Value = "Cacotemnus thomsoni"   # G61
$ws.Cells.Item(61, 8).Value = "(Kraatz, 1881)"        # H61
$ws.Cells.Item(61, 17).Value = 571834                # Q61
$ws.Cells.Item(61, 18).Value = 6697641               # R61

# --- Row 62 ---------------------------------------------------------
$ws.Cells.Item(62, 1).Value = 112230613              # A62
$ws.Cells.Item(62, 2).Value = 89553                  # B62
$ws.Cells.Item(62, 4).Value = "NT"                    # D62
$ws.Cells.Item(62, 5).Value = 1202                   # E62
$ws.Cells.Item(62, 6).Value = "Ullticka"              # F62
$ws.Cells.Item(62, 7).Value = "Phellinidium ferrugineofuscum" # G62
$ws.Cells.Item(62, 8).Value = "(P.Karst.) Fiasson & Niemelä"  # H62
$ws.Cells.Item(62, 17).Value = 571799                # Q62
$ws.Cells.Item(62, 18).Value = 6697620               # R62
